$d = $word.ActiveDocument

# --- 1. "This requires Visual Studio 2012 Express for Desktop" -> 2013 ---
$d.Content.Find.Execute("Visual Studio 2012 Express for Desktop", $true, $false, $false, $false, $false, $true, 1, $false, "Visual Studio 2013 Express for Desktop", 2)

# --- 2. "wxWidgets 2.9 SVN, 32 bit." -> "wxWidgets 3.0 SVN, 32 bit." ---
$d.Content.Find.Execute(" 2.9 SVN, 32 bit.", $true, $false, $false, $false, $false, $true, 1, $false, " 3.0 SVN, 32 bit.", 2)

# --- 3. Delete whole paragraphs that are dropped, from bottom to top to keep indices stable ---
# 17: Copy the setup0.h file to setup.h
$d.Paragraphs(17).Range.Delete()
# 16: Open the folder C:\wxWidgets32\include\wx\msw in Windows Explorer
$d.Paragraphs(16).Range.Delete()
# 13: Open Visual Studio 2012, and open the c:\wxWidgets32\build\msw\wx_vc9.sln file
$d.Paragraphs(13).Range.Delete()
# 12: When wxWidgets 2.9.5 comes out ...
$d.Paragraphs(12).Range.Delete()
# 11: I checked out revision 73160. ...
$d.Paragraphs(11).Range.Delete()
# 10: It'll take some time to get all the files
$d.Paragraphs(10).Range.Delete()
# 7: Create a folder c:\wxWidgets32
$d.Paragraphs(7).Range.Delete()

Write-Output "deleted"

# --- 4. "Right click on it, and check out (via TortoiseSVN) the wxWidgets trunk from the public SVN"
#        -> "Download wxWidgets 3.0 from "
$d.Content.Find.Execute("Right click on it, and check out (via TortoiseSVN) the ", $true, $false, $false, $false, $false, $true, 1, $false, "Download ", 2)
$d.Content.Find.Execute(" trunk from the public SVN", $true, $false, $false, $false, $false, $true, 1, $false, " 3.0 from ", 2)

# --- 5. hyperlink text + remove list numbering on that paragraph ---
$d.Content.Find.Execute("http://svn.wxwidgets.org/svn/wx/wxWidgets/trunk", $true, $false, $false, $false, $false, $true, 1, $false, "https://sourceforge.net/projects/wxwindows/files/3.0.0/wxWidgets-3.0.0.zip/download", 2)
$d.Paragraphs(8).Range.ListFormat.RemoveNumbers()

Write-Output "rewrote links"

# --- 6. Insert two new steps after the hyperlink paragraph (currently paragraph 8):
#        "Extract to wxWidgets-3.0.0-x32" and
#        "Open Visual Studio 2013, and open the C:\wxWidgets-3.0.0-x32\build\msw file"
$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphBefore()
$newp1 = $d.Paragraphs(9)
$newp1.Range.ListFormat.ListLevelNumber = 1
$newp1.Range.Text = "Extract to wxWidgets-3.0.0-x32"

$d.Paragraphs(10).Range.InsertParagraphBefore()
$newp2 = $d.Paragraphs(10)
$newp2.Range.ListFormat.ListLevelNumber = 1
$newp2.Range.Text = "Open Visual Studio 2013, and open the C:\wxWidgets-3.0.0-x32\build\msw file"

Write-Output "inserted steps"

# --- 7. " Set the WXMSW3 environment variable to c:\wxWidgets32" -> "...to C:\wxWidgets-3.0.0-x32"
$d.Content.Find.Execute("Set the WXMSW3 environment variable to c:\wxWidgets32", $true, $false, $false, $false, $false, $true, 1, $false, "Set the WXMSW3 environment variable to C:\wxWidgets-3.0.0-x32", 2)

Write-Output "wxmsw3 rewritten"


